$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$arr = New-Object 'object[,]' 64,3
$arr[0,0] = '(''pretty'', ''woman'')'
$arr[0,1] = 14
$arr[0,2] = 0.0396600566572238
$arr[1,0] = '(''low'', ''budget'')'
$arr[1,1] = 12
$arr[1,2] = 0.0339943342776204
$arr[2,0] = '(''waste'', ''time'')'
$arr[2,1] = 12
$arr[2,2] = 0.0339943342776204
$arr[3,0] = '(''bad'', ''movie'')'
$arr[3,1] = 10
$arr[3,2] = 0.028328611898017
$arr[4,0] = '(''mst'', ''k'')'
$arr[4,1] = 10
$arr[4,2] = 0.028328611898017
$arr[5,0] = '(''whole'', ''thing'')'
$arr[5,1] = 10
$arr[5,2] = 0.028328611898017
$arr[6,0] = '(''camera'', ''work'')'
$arr[6,1] = 9
$arr[6,2] = 0.0254957507082153
$arr[7,0] = '(''movie'', ''ever'')'
$arr[7,1] = 9
$arr[7,2] = 0.0254957507082153
$arr[8,0] = '(''dick'', ''jane'')'
$arr[8,1] = 8
$arr[8,2] = 0.0226628895184136
$arr[9,0] = '(''pretty'', ''much'')'
$arr[9,1] = 8
$arr[9,2] = 0.0226628895184136
$arr[10,0] = '(''movie'', ''would'')'
$arr[10,1] = 7
$arr[10,2] = 0.0198300283286119
$arr[11,0] = '(''really'', ''really'')'
$arr[11,1] = 7
$arr[11,2] = 0.0198300283286119
$arr[12,0] = '(''main'', ''character'')'
$arr[12,1] = 6
$arr[12,2] = 0.0169971671388102
$arr[13,0] = '(''r'', ''rated'')'
$arr[13,1] = 6
$arr[13,2] = 0.0169971671388102
$arr[14,0] = '(''look'', ''like'')'
$arr[14,1] = 6
$arr[14,2] = 0.0169971671388102
$arr[15,0] = '(''really'', ''bad'')'
$arr[15,1] = 6
$arr[15,2] = 0.0169971671388102
$arr[16,0] = '(''one'', ''would'')'
$arr[16,1] = 6
$arr[16,2] = 0.0169971671388102
$arr[17,0] = '(''great'', ''movie'')'
$arr[17,1] = 5
$arr[17,2] = 0.0141643059490085
$arr[18,0] = '(''horror'', ''film'')'
$arr[18,1] = 5
$arr[18,2] = 0.0141643059490085
$arr[19,0] = '(''felt'', ''like'')'
$arr[19,1] = 5
$arr[19,2] = 0.0141643059490085
$arr[20,0] = '(''film'', ''like'')'
$arr[20,1] = 5
$arr[20,2] = 0.0141643059490085
$arr[21,0] = '(''first'', ''one'')'
$arr[21,1] = 5
$arr[21,2] = 0.0141643059490085
$arr[22,0] = '(''movie'', ''bad'')'
$arr[22,1] = 5
$arr[22,2] = 0.0141643059490085
$arr[23,0] = '(''worst'', ''movie'')'
$arr[23,1] = 5
$arr[23,2] = 0.0141643059490085
$arr[24,0] = '(''character'', ''development'')'
$arr[24,1] = 5
$arr[24,2] = 0.0141643059490085
$arr[25,0] = '(''rishi'', ''kapoor'')'
$arr[25,1] = 5
$arr[25,2] = 0.0141643059490085
$arr[26,0] = '(''film'', ''without'')'
$arr[26,1] = 5
$arr[26,2] = 0.0141643059490085
$arr[27,0] = '(''first'', ''half'')'
$arr[27,1] = 5
$arr[27,2] = 0.0141643059490085
$arr[28,0] = '(''waste'', ''money'')'
$arr[28,1] = 5
$arr[28,2] = 0.0141643059490085
$arr[29,0] = '(''sex'', ''scenes'')'
$arr[29,1] = 5
$arr[29,2] = 0.0141643059490085
$arr[30,0] = '(''movie'', ''still'')'
$arr[30,1] = 5
$arr[30,2] = 0.0141643059490085
$arr[31,0] = '(''whole'', ''movie'')'
$arr[31,1] = 5
$arr[31,2] = 0.0141643059490085
$arr[32,0] = '(''movie'', ''even'')'
$arr[32,1] = 5
$arr[32,2] = 0.0141643059490085
$arr[33,0] = '(''happen'', ''film'')'
$arr[33,1] = 5
$arr[33,2] = 0.0141643059490085
$arr[34,0] = '(''see'', ''dick'')'
$arr[34,1] = 5
$arr[34,2] = 0.0141643059490085
$arr[35,0] = '(''shahrukh'', ''khan'')'
$arr[35,1] = 5
$arr[35,2] = 0.0141643059490085
$arr[36,0] = '(''bu'', ''uel'')'
$arr[36,1] = 4
$arr[36,2] = 0.0113314447592068
$arr[37,0] = '(''never'', ''get'')'
$arr[37,1] = 4
$arr[37,2] = 0.0113314447592068
$arr[38,0] = '(''first'', ''movie'')'
$arr[38,1] = 4
$arr[38,2] = 0.0113314447592068
$arr[39,0] = '(''want'', ''watch'')'
$arr[39,1] = 4
$arr[39,2] = 0.0113314447592068
$arr[40,0] = '(''pretty'', ''bad'')'
$arr[40,1] = 4
$arr[40,2] = 0.0113314447592068
$arr[41,0] = '(''half'', ''hour'')'
$arr[41,1] = 4
$arr[41,2] = 0.0113314447592068
$arr[42,0] = '(''minutes'', ''movie'')'
$arr[42,1] = 4
$arr[42,2] = 0.0113314447592068
$arr[43,0] = '(''looked'', ''like'')'
$arr[43,1] = 4
$arr[43,2] = 0.0113314447592068
$arr[44,0] = '(''yet'', ''another'')'
$arr[44,1] = 4
$arr[44,2] = 0.0113314447592068
$arr[45,0] = '(''hedy'', ''lamarr'')'
$arr[45,1] = 4
$arr[45,2] = 0.0113314447592068
$arr[46,0] = '(''could'', ''good'')'
$arr[46,1] = 4
$arr[46,2] = 0.0113314447592068
$arr[47,0] = '(''still'', ''think'')'
$arr[47,1] = 4
$arr[47,2] = 0.0113314447592068
$arr[48,0] = '(''science'', ''fiction'')'
$arr[48,1] = 4
$arr[48,2] = 0.0113314447592068
$arr[49,0] = '(''awful'', ''movie'')'
$arr[49,1] = 4
$arr[49,2] = 0.0113314447592068
$arr[50,0] = '(''say'', ''film'')'
$arr[50,1] = 4
$arr[50,2] = 0.0113314447592068
$arr[51,0] = '(''original'', ''film'')'
$arr[51,1] = 4
$arr[51,2] = 0.0113314447592068
$arr[52,0] = '(''still'', ''like'')'
$arr[52,1] = 4
$arr[52,2] = 0.0113314447592068
$arr[53,0] = '(''good'', ''enough'')'
$arr[53,1] = 4
$arr[53,2] = 0.0113314447592068
$arr[54,0] = '(''like'', ''could'')'
$arr[54,1] = 4
$arr[54,2] = 0.0113314447592068
$arr[55,0] = '(''year'', ''olds'')'
$arr[55,1] = 4
$arr[55,2] = 0.0113314447592068
$arr[56,0] = '(''every'', ''time'')'
$arr[56,1] = 4
$arr[56,2] = 0.0113314447592068
$arr[57,0] = '(''say'', ''movie'')'
$arr[57,1] = 4
$arr[57,2] = 0.0113314447592068
$arr[58,0] = '(''mary'', ''lou'')'
$arr[58,1] = 4
$arr[58,2] = 0.0113314447592068
$arr[59,0] = '(''bad'', ''acting'')'
$arr[59,1] = 4
$arr[59,2] = 0.0113314447592068
$arr[60,0] = '(''movie'', ''almost'')'
$arr[60,1] = 4
$arr[60,2] = 0.0113314447592068
$arr[61,0] = '(''many'', ''people'')'
$arr[61,1] = 4
$arr[61,2] = 0.0113314447592068
$arr[62,0] = '(''watching'', ''movie'')'
$arr[62,1] = 4
$arr[62,2] = 0.0113314447592068
$arr[63,0] = '(''really'', ''think'')'
$arr[63,1] = 4
$arr[63,2] = 0.0113314447592068
$ws.Range("A2:C65").Value = $arr

$ws = $wb.Worksheets.Item(2)
$arr = New-Object 'object[,]' 64,3
$arr[0,0] = '(''tony'', ''hawk'')'
$arr[0,1] = 35
$arr[0,2] = 0.0831353919239905
$arr[1,0] = '(''pro'', ''skater'')'
$arr[1,1] = 22
$arr[1,2] = 0.05225653206650831
$arr[2,0] = '(''hawk'', ''pro'')'
$arr[2,1] = 22
$arr[2,2] = 0.05225653206650831
$arr[3,0] = '(''skater'', ''x'')'
$arr[3,1] = 19
$arr[3,2] = 0.04513064133016627
$arr[4,0] = '(''one'', ''best'')'
$arr[4,1] = 17
$arr[4,2] = 0.04038004750593824
$arr[5,0] = '(''one'', ''thing'')'
$arr[5,1] = 11
$arr[5,2] = 0.02612826603325416
$arr[6,0] = '(''see'', ''movie'')'
$arr[6,1] = 9
$arr[6,2] = 0.02137767220902613
$arr[7,0] = '(''first'', ''time'')'
$arr[7,1] = 8
$arr[7,2] = 0.01900237529691211
$arr[8,0] = '(''james'', ''bond'')'
$arr[8,1] = 7
$arr[8,2] = 0.0166270783847981
$arr[9,0] = '(''great'', ''job'')'
$arr[9,1] = 7
$arr[9,2] = 0.0166270783847981
$arr[10,0] = '(''first'', ''film'')'
$arr[10,1] = 7
$arr[10,2] = 0.0166270783847981
$arr[11,0] = '(''say'', ''never'')'
$arr[11,1] = 6
$arr[11,2] = 0.01425178147268409
$arr[12,0] = '(''highly'', ''recommend'')'
$arr[12,1] = 6
$arr[12,2] = 0.01425178147268409
$arr[13,0] = '(''chipmunk'', ''adventure'')'
$arr[13,1] = 6
$arr[13,2] = 0.01425178147268409
$arr[14,0] = '(''never'', ''say'')'
$arr[14,1] = 6
$arr[14,2] = 0.01425178147268409
$arr[15,0] = '(''movie'', ''made'')'
$arr[15,1] = 6
$arr[15,2] = 0.01425178147268409
$arr[16,0] = '(''x'', ''career'')'
$arr[16,1] = 6
$arr[16,2] = 0.01425178147268409
$arr[17,0] = '(''black'', ''white'')'
$arr[17,1] = 6
$arr[17,2] = 0.01425178147268409
$arr[18,0] = '(''ghost'', ''stories'')'
$arr[18,1] = 6
$arr[18,2] = 0.01425178147268409
$arr[19,0] = '(''get'', ''see'')'
$arr[19,1] = 6
$arr[19,2] = 0.01425178147268409
$arr[20,0] = '(''first'', ''saw'')'
$arr[20,1] = 6
$arr[20,2] = 0.01425178147268409
$arr[21,0] = '(''takes'', ''place'')'
$arr[21,1] = 5
$arr[21,2] = 0.01187648456057007
$arr[22,0] = '(''movie'', ''good'')'
$arr[22,1] = 5
$arr[22,2] = 0.01187648456057007
$arr[23,0] = '(''like'', ''one'')'
$arr[23,1] = 5
$arr[23,2] = 0.01187648456057007
$arr[24,0] = '(''well'', ''made'')'
$arr[24,1] = 5
$arr[24,2] = 0.01187648456057007
$arr[25,0] = '(''subject'', ''matter'')'
$arr[25,1] = 5
$arr[25,2] = 0.01187648456057007
$arr[26,0] = '(''films'', ''ever'')'
$arr[26,1] = 5
$arr[26,2] = 0.01187648456057007
$arr[27,0] = '(''never'', ''seen'')'
$arr[27,1] = 5
$arr[27,2] = 0.01187648456057007
$arr[28,0] = '(''feel'', ''like'')'
$arr[28,1] = 5
$arr[28,2] = 0.01187648456057007
$arr[29,0] = '(''bobby'', ''earl'')'
$arr[29,1] = 5
$arr[29,2] = 0.01187648456057007
$arr[30,0] = '(''main'', ''characters'')'
$arr[30,1] = 5
$arr[30,2] = 0.01187648456057007
$arr[31,0] = '(''never'', ''heard'')'
$arr[31,1] = 5
$arr[31,2] = 0.01187648456057007
$arr[32,0] = '(''old'', ''man'')'
$arr[32,1] = 5
$arr[32,2] = 0.01187648456057007
$arr[33,0] = '(''world'', ''war'')'
$arr[33,1] = 5
$arr[33,2] = 0.01187648456057007
$arr[34,0] = '(''anything'', ''else'')'
$arr[34,1] = 5
$arr[34,2] = 0.01187648456057007
$arr[35,0] = '(''film'', ''made'')'
$arr[35,1] = 5
$arr[35,2] = 0.01187648456057007
$arr[36,0] = '(''game'', ''play'')'
$arr[36,1] = 5
$arr[36,2] = 0.01187648456057007
$arr[37,0] = '(''movies'', ''ever'')'
$arr[37,1] = 5
$arr[37,2] = 0.01187648456057007
$arr[38,0] = '(''hawk'', ''career'')'
$arr[38,1] = 5
$arr[38,2] = 0.01187648456057007
$arr[39,0] = '(''one'', ''film'')'
$arr[39,1] = 5
$arr[39,2] = 0.01187648456057007
$arr[40,0] = '(''best'', ''friend'')'
$arr[40,1] = 5
$arr[40,2] = 0.01187648456057007
$arr[41,0] = '(''must'', ''see'')'
$arr[41,1] = 5
$arr[41,2] = 0.01187648456057007
$arr[42,0] = '(''gene'', ''kelly'')'
$arr[42,1] = 5
$arr[42,2] = 0.01187648456057007
$arr[43,0] = '(''one'', ''night'')'
$arr[43,1] = 5
$arr[43,2] = 0.01187648456057007
$arr[44,0] = '(''star'', ''trek'')'
$arr[44,1] = 5
$arr[44,2] = 0.01187648456057007
$arr[45,0] = '(''second'', ''time'')'
$arr[45,1] = 5
$arr[45,2] = 0.01187648456057007
$arr[46,0] = '(''like'', ''really'')'
$arr[46,1] = 4
$arr[46,2] = 0.009501187648456057
$arr[47,0] = '(''phantom'', ''lady'')'
$arr[47,1] = 4
$arr[47,2] = 0.009501187648456057
$arr[48,0] = '(''seen'', ''movie'')'
$arr[48,1] = 4
$arr[48,2] = 0.009501187648456057
$arr[49,0] = '(''like'', ''see'')'
$arr[49,1] = 4
$arr[49,2] = 0.009501187648456057
$arr[50,0] = '(''beginning'', ''end'')'
$arr[50,1] = 4
$arr[50,2] = 0.009501187648456057
$arr[51,0] = '(''wagon'', ''master'')'
$arr[51,1] = 4
$arr[51,2] = 0.009501187648456057
$arr[52,0] = '(''really'', ''enjoyed'')'
$arr[52,1] = 4
$arr[52,2] = 0.009501187648456057
$arr[53,0] = '(''well'', ''worth'')'
$arr[53,1] = 4
$arr[53,2] = 0.009501187648456057
$arr[54,0] = '(''great'', ''work'')'
$arr[54,1] = 4
$arr[54,2] = 0.009501187648456057
$arr[55,0] = '(''young'', ''woman'')'
$arr[55,1] = 4
$arr[55,2] = 0.009501187648456057
$arr[56,0] = '(''old'', ''fashioned'')'
$arr[56,1] = 4
$arr[56,2] = 0.009501187648456057
$arr[57,0] = '(''action'', ''films'')'
$arr[57,1] = 4
$arr[57,2] = 0.009501187648456057
$arr[58,0] = '(''least'', ''one'')'
$arr[58,1] = 4
$arr[58,2] = 0.009501187648456057
$arr[59,0] = '(''gave'', ''us'')'
$arr[59,1] = 4
$arr[59,2] = 0.009501187648456057
$arr[60,0] = '(''works'', ''well'')'
$arr[60,1] = 4
$arr[60,2] = 0.009501187648456057
$arr[61,0] = '(''jack'', ''frost'')'
$arr[61,1] = 4
$arr[61,2] = 0.009501187648456057
$arr[62,0] = '(''bad'', ''guys'')'
$arr[62,1] = 4
$arr[62,2] = 0.009501187648456057
$arr[63,0] = '(''wife'', ''daughter'')'
$arr[63,1] = 4
$arr[63,2] = 0.009501187648456057
$ws.Range("A2:C65").Value = $arr

Write-Output "done"
